$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 4750
$ws.Range("J43").Value = 4750
$ws.Range("L43").Value = 4750
$ws.Range("N43").Value = -4888
$ws.Range("H48").Value = 2787.7
$ws.Range("J48").Value = 2787.7
$ws.Range("L48").Value = 8363.099999999999
$ws.Range("N48").Value = -8947.099999999999
$ws.Range("H56").Value = 2787.7
$ws.Range("J56").Value = 2787.7
$ws.Range("L56").Value = 8363.099999999999
$ws.Range("N56").Value = -9431.099999999999
$ws.Range("H62").Value = 72998.3
$ws.Range("I62").Value = 72998.3
$ws.Range("K62").Value = 72998.3
$ws.Range("M62").Value = -72374.3
$ws.Range("H65").Value = 72998.3
$ws.Range("I65").Value = 72998.3
$ws.Range("K65").Value = 364991.5
$ws.Range("M65").Value = -361871.5
$ws.Range("H80").Value = 150849.75
$ws.Range("I80").Value = 150849.75
$ws.Range("K80").Value = 452549.25
$ws.Range("M80").Value = -451551.25
$ws.Range("H83").Value = 150849.75
$ws.Range("I83").Value = 150849.75
$ws.Range("K83").Value = 1357647.75
$ws.Range("M83").Value = -1352655.75
$ws.Range("H94").Value = 6285.2856
$ws.Range("I94").Value = 4799.4
$ws.Range("J94").Value = 10000
$ws.Range("K94").Value = 4799.4
$ws.Range("L94").Value = 10000
$ws.Range("M94").Value = -4348.4
$ws.Range("N94").Value = -10902
$ws.Range("H129").Value = 5295.2607
$ws.Range("J129").Value = 10936.111
$ws.Range("L129").Value = 32808.333
$ws.Range("N129").Value = -42808.333
$ws.Range("H131").Value = 2798.6843
$ws.Range("I131").Value = 1261.3125
$ws.Range("J131").Value = 10998
$ws.Range("K131").Value = 3783.9375
$ws.Range("L131").Value = 32994
$ws.Range("M131").Value = 1256.0625
$ws.Range("N131").Value = -43074
$ws.Range("H132").Value = 6924.1226
$ws.Range("I132").Value = 5132.5
$ws.Range("K132").Value = 15397.5
$ws.Range("M132").Value = -12867.5
$ws.Range("H137").Value = 3769.9333
$ws.Range("I137").Value = 3526
$ws.Range("K137").Value = 10578
$ws.Range("M137").Value = -8028
$ws.Range("H138").Value = 7171.914
$ws.Range("I138").Value = 8645.714
$ws.Range("J138").Value = 6803.4644
$ws.Range("K138").Value = 25937.142
$ws.Range("L138").Value = 20410.3932
$ws.Range("M138").Value = -20797.142
$ws.Range("N138").Value = -30690.3932

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 32040.5
$ws.Range("I32").Value = 21123.75
$ws.Range("J32").Value = 53874
$ws.Range("K32").Value = 21123.75
$ws.Range("L32").Value = 53874
$ws.Range("M32").Value = -20836.75
$ws.Range("N32").Value = -54448
$ws.Range("H88").Value = 2075.8
$ws.Range("I88").Value = 1489
$ws.Range("J88").Value = 2222.5
$ws.Range("K88").Value = 1489
$ws.Range("L88").Value = 2222.5
$ws.Range("M88").Value = -1083
$ws.Range("N88").Value = -3034.5
$ws.Range("H91").Value = 2075.8
$ws.Range("I91").Value = 1489
$ws.Range("J91").Value = 2222.5
$ws.Range("K91").Value = 1489
$ws.Range("L91").Value = 2222.5
$ws.Range("M91").Value = -85
$ws.Range("N91").Value = -5030.5
$ws.Range("H110").Value = 4400.8
$ws.Range("I110").Value = 5152
$ws.Range("K110").Value = 5152
$ws.Range("M110").Value = -3107
$ws.Range("H122").Value = 3610.7778
$ws.Range("I122").Value = 3610.7778
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 10832.3334
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -8382.3334
$ws.Range("N122").ClearContents()
$ws.Range("H132").Value = 28619.77
$ws.Range("I132").Value = 30540.695
$ws.Range("K132").Value = 91622.08499999999
$ws.Range("M132").Value = -89092.08499999999

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H74").Value = 34998.5
$ws.Range("J74").Value = 34998.5
$ws.Range("L74").Value = 34998.5
$ws.Range("N74").Value = -36870.5
$ws.Range("H77").Value = 34998.5
$ws.Range("J77").Value = 34998.5
$ws.Range("L77").Value = 104995.5
$ws.Range("N77").Value = -114355.5
$ws.Range("H86").Value = 2617.25
$ws.Range("I86").Value = 1905.7778
$ws.Range("K86").Value = 1905.7778
$ws.Range("M86").Value = -782.7778000000001
$ws.Range("H89").Value = 2617.25
$ws.Range("I89").Value = 1905.7778
$ws.Range("K89").Value = 9528.889000000001
$ws.Range("M89").Value = -3912.889000000001
$ws.Range("H134").Value = 2463.3333
$ws.Range("I134").Value = 2463.3333
$ws.Range("K134").Value = 7389.999899999999
$ws.Range("M134").Value = -4854.999899999999

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2908.611
$ws.Range("I31").Value = 2508.2144
$ws.Range("K31").Value = 2508.2144
$ws.Range("M31").Value = -2213.2144
$ws.Range("H34").Value = 2908.611
$ws.Range("I34").Value = 2508.2144
$ws.Range("K34").Value = 2508.2144
$ws.Range("M34").Value = -2306.2144
$ws.Range("H105").Value = 2999.75
$ws.Range("J105").Value = 3000
$ws.Range("L105").Value = 3000
$ws.Range("N105").Value = -6494
$ws.Range("H131").Value = 41158.4
$ws.Range("J131").Value = 41374.25
$ws.Range("L131").Value = 41374.25
$ws.Range("N131").Value = -51454.25
$ws.Range("H134").Value = 35406.094
$ws.Range("I134").Value = 42353.848
$ws.Range("K134").Value = 127061.544
$ws.Range("M134").Value = -124526.544
$ws.Range("H139").Value = 79899.5
$ws.Range("I139").Value = 79800
$ws.Range("K139").Value = 79800
$ws.Range("M139").Value = -74660
$ws.Range("H141").Value = 549812.9399999999
$ws.Range("J141").Value = 576024.75
$ws.Range("L141").Value = 576024.75
$ws.Range("N141").Value = -586384.75

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 43.0625
$ws.Range("J38").Value = 38.57143
$ws.Range("L38").Value = 115.71429
$ws.Range("N38").Value = -809.71429
$ws.Range("H107").Value = 506.5
$ws.Range("I107").Value = 157.5
$ws.Range("J107").Value = 681
$ws.Range("K107").Value = 472.5
$ws.Range("L107").Value = 2043
$ws.Range("M107").Value = 1447.5
$ws.Range("N107").Value = -5883
$ws.Range("H114").Value = 1455.375
$ws.Range("I114").Value = 581.5
$ws.Range("J114").Value = 1979.7
$ws.Range("K114").Value = 1744.5
$ws.Range("L114").Value = 5939.1
$ws.Range("M114").Value = 1509.5
$ws.Range("N114").Value = -12447.1
$ws.Range("H119").Value = 8329
$ws.Range("I119").Value = 6384
$ws.Range("K119").Value = 19152
$ws.Range("M119").Value = -14314
$ws.Range("H120").Value = 15172.875
$ws.Range("J120").Value = 14304.111
$ws.Range("L120").Value = 42912.333
$ws.Range("N120").Value = -52588.333
$ws.Range("H121").Value = 113062
$ws.Range("I121").Value = 3676.6667
$ws.Range("J121").Value = 167754.67
$ws.Range("K121").Value = 11030.0001
$ws.Range("L121").Value = 503264.01
$ws.Range("M121").Value = -9720.000100000001
$ws.Range("N121").Value = -505884.01
$ws.Range("H129").Value = 1215945.4
$ws.Range("J129").Value = 8500000
$ws.Range("L129").Value = 25500000
$ws.Range("N129").Value = -25510000
$ws.Range("H131").Value = 4010618.5
$ws.Range("I131").Value = 1858
$ws.Range("J131").Value = 5012808.5
$ws.Range("K131").Value = 5574
$ws.Range("L131").Value = 15038425.5
$ws.Range("M131").Value = -534
$ws.Range("N131").Value = -15048505.5
$ws.Range("H137").Value = 1999.5
$ws.Range("I137").Value = 2000
$ws.Range("J137").Value = 1999
$ws.Range("K137").Value = 6000
$ws.Range("L137").Value = 5997
$ws.Range("M137").Value = -900
$ws.Range("N137").Value = -16197

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3776.2307
$ws.Range("I122").Value = 3553.7273
$ws.Range("K122").Value = 10661.1819
$ws.Range("M122").Value = -8211.1819
$ws.Range("H132").Value = 94949.5
$ws.Range("I132").Value = 102899.45
$ws.Range("J132").Value = 7500
$ws.Range("K132").Value = 308698.35
$ws.Range("L132").Value = 22500
$ws.Range("M132").Value = -306168.35
$ws.Range("N132").Value = -27560

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 4214.967
$ws.Range("I22").Value = 3254.0908
$ws.Range("J22").Value = 4771.263
$ws.Range("K22").Value = 3254.0908
$ws.Range("L22").Value = 4771.263
$ws.Range("M22").Value = -2959.0908
$ws.Range("N22").Value = -5361.263
$ws.Range("H27").Value = 4214.967
$ws.Range("I27").Value = 3254.0908
$ws.Range("J27").Value = 4771.263
$ws.Range("K27").Value = 3254.0908
$ws.Range("L27").Value = 4771.263
$ws.Range("M27").Value = -3147.0908
$ws.Range("N27").Value = -4985.263

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 19999
$ws.Range("I51").Value = 19999
$ws.Range("K51").Value = 19999
$ws.Range("M51").Value = -19489
$ws.Range("H122").Value = 3234.2778
$ws.Range("I122").Value = 2663.5625
$ws.Range("K122").Value = 7990.6875
$ws.Range("M122").Value = -5540.6875
$ws.Range("H132").Value = 37164.51
$ws.Range("I132").Value = 39180.098
$ws.Range("K132").Value = 117540.294
$ws.Range("M132").Value = -115010.294
